# Add 2022-Q4 data:
#  - The existing "2022-Q3" sheet becomes "2022-Q4" and gets the new quarter's data.
#  - A fresh "2022-Q3" sheet is inserted right after it, preserving the old data.
#  - The "总计" (totals) sheet gets a new row for 2022-Q3 and its 2022-Q4 row is refreshed.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ = $wb.Worksheets.Item(2)

# 1) Insert a brand-new worksheet right after the existing quarter sheet; this
#    will hold the previous (2022-Q3) data and becomes the active sheet, just
#    like Excel does after Worksheets.Add.
$wsQ3 = $wb.Worksheets.Add($null, $wsQ)

# 2) Copy the old sheet's full contents + formatting into the new sheet before
#    the old sheet's data is overwritten with the new quarter's numbers, then
#    drop the stray blank A1 anchor cell the copy leaves behind (the source
#    sheet's data starts at column B).
$wsQ.UsedRange.Copy($wsQ3.Range("A1"))
$wsQ3.Range("A1").ClearContents()

# 3) Rename the original quarter sheet to the new quarter first, freeing up
#    the "2022-Q3" name for the new sheet (renaming the new sheet first would
#    collide with the still-named-"2022-Q3" original sheet).
$wsQ.Name = "2022-Q4"
$wsQ3.Name = "2022-Q3"

# 4) Overwrite the (now 2022-Q4) sheet with the new holdings data. Columns
#    B and D:G hold numeric-looking text (fund codes keep leading zeros,
#    percentages keep trailing zeros) rather than numbers, so force text
#    formatting before writing so they aren't silently coerced to numbers
#    (column C is always non-numeric text, so it needs no such nudge).
# (applied as two single-area ranges - NumberFormat on a multi-area Range
# only takes effect on its first area.)
$wsQ.Range("B2:B4").NumberFormat = "@"
$wsQ.Range("D2:G4").NumberFormat = "@"

$wsQ.Range("B2").Value = "003835"
$wsQ.Range("C2").Value = "鹏华沪深港新兴成长灵活配置混合"
$wsQ.Range("D2").Value = "45.58"
$wsQ.Range("E2").Value = "94.28"
$wsQ.Range("F2").Value = "4.08"
$wsQ.Range("G2").Value = "1.8597"
$wsQ.Range("H2").Value = 5

$wsQ.Range("B3").Value = "013250"
$wsQ.Range("C3").Value = "红土创新智能制造混合"
$wsQ.Range("D3").Value = "1.19"
$wsQ.Range("E3").Value = "90.41"
$wsQ.Range("F3").Value = "5.50"
$wsQ.Range("G3").Value = "0.0654"
$wsQ.Range("H3").Value = 5

$wsQ.Range("B4").Value = "004044"
$wsQ.Range("C4").Value = "金鹰转型动力灵活配置混合"
$wsQ.Range("D4").Value = "0.65"
$wsQ.Range("E4").Value = "90.51"
$wsQ.Range("F4").Value = "7.29"
$wsQ.Range("G4").Value = "0.0474"
$wsQ.Range("H4").Value = 1

# 4b) The "@" text format above leaves those cells pointing at a dedicated
#     "Text" style, whereas the original data cells carry no explicit style
#     at all. Paste the plain (unstyled) format from a guaranteed-blank cell
#     back over them - this only touches formatting, so the text already in
#     the cells (with leading/trailing zeros intact) is left alone.
$wsTotal.Range("Z100").Copy()
$wsQ.Range("B2:B4").PasteSpecial(-4122)
$wsQ.Range("D2:G4").PasteSpecial(-4122)
$wsTotal.Range("Z100").ClearContents()

# 5) Match the "总计" sheet's header/first-column formatting on the 2022-Q4
#    sheet (it previously used the other style since it was built fresh).
$wsTotal.Range("B1").Copy()
$wsQ.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ.Range("A2:A4").PasteSpecial(-4122)

# 6) Update the totals sheet: refresh the 2022-Q4 row and append a 2022-Q3 row.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 1.97

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 0.22
